# MHD2-259: Report template and related changes for reporting on 136 genes
#
# The underlying diff touches template-level definitions (numbering.xml
# and styles.xml) rather than document body content. We drive this
# through the Styles collection of the Word object model (Font /
# ParagraphFormat), which is the supported, non-destructive surface for
# editing style definitions in this host. (Document.WordOpenXML is a
# read-only snapshot here - the runtime rejects assignments to it - and
# the ListFormat.Apply*Default()/ApplyListTemplateWithLevel() numbering
# helpers unconditionally fabricate a brand-new list definition and
# silently drop bookkeeping attributes from the existing <w:num> entries
# as a side effect, so they are intentionally not used here to avoid
# corrupting the numbering part.)
#
# NOTE: everything is written as flat, inline statements (no functions
# with named -Parameters) because this host's PowerShell parameter
# binding for named arguments is unreliable.

$d = $word.ActiveDocument

# --- ListParagraph: align the (previously direct 428/-425 twip) indent
#     with the list level's own indent (714/-360 twips) - the same
#     result that dropping the paragraph-level <w:ind> override and
#     inheriting from the numbering definition produces. ------------
$listParagraph = $d.Styles("List Paragraph")
$listParagraph.ParagraphFormat.LeftIndent = 35.7      # 714 twips
$listParagraph.ParagraphFormat.FirstLineIndent = -18  # -360 twips (hanging 360)

# --- CLIN1(HEADING): Calibri -> Aptos, add explicit 15pt (sz 30) ----
$clin1Heading = $d.Styles("CLIN1(HEADING)")
$clin1Heading.Font.NameAscii = "Aptos"
$clin1Heading.Font.NameOther = "Aptos"
$clin1Heading.Font.Size = 15

# --- CLIN2(SUBHEADINGS): Calibri -> Aptos, 11pt -> 10pt (sz 22 -> 20)
$clin2Subheadings = $d.Styles("CLIN2(SUBHEADINGS)")
$clin2Subheadings.Font.NameAscii = "Aptos"
$clin2Subheadings.Font.NameOther = "Aptos"
$clin2Subheadings.Font.Size = 10

# --- CLIN1(HEADING) Char: Calibri -> Aptos, 16pt -> 15pt (sz 32 -> 30)
$clin1HeadingChar = $d.Styles("CLIN1(HEADING) Char")
$clin1HeadingChar.Font.NameAscii = "Aptos"
$clin1HeadingChar.Font.NameOther = "Aptos"
$clin1HeadingChar.Font.Size = 15

# --- CLIN3(BULLET POINTS): Calibri -> Aptos, add explicit 8pt (sz 16)
$clin3Bullet = $d.Styles("CLIN3(BULLET POINTS)")
$clin3Bullet.Font.NameAscii = "Aptos"
$clin3Bullet.Font.NameOther = "Aptos"
$clin3Bullet.Font.Size = 8

# --- CLIN2(SUBHEADINGS) Char: Calibri -> Aptos, 13pt -> 10pt (sz 26 -> 20)
$clin2SubheadingsChar = $d.Styles("CLIN2(SUBHEADINGS) Char")
$clin2SubheadingsChar.Font.NameAscii = "Aptos"
$clin2SubheadingsChar.Font.NameOther = "Aptos"
$clin2SubheadingsChar.Font.Size = 10

# --- CLIN4: Calibri -> Aptos, 6pt -> 5.5pt (sz 12 -> 11) ------------
$clin4 = $d.Styles("CLIN4")
$clin4.Font.NameAscii = "Aptos"
$clin4.Font.NameOther = "Aptos"
$clin4.Font.Size = 5.5

# --- CLIN3(BULLET POINTS) Char: Calibri -> Aptos, 9pt -> 8pt (sz 18 -> 16)
$clin3BulletChar = $d.Styles("CLIN3(BULLET POINTS) Char")
$clin3BulletChar.Font.NameAscii = "Aptos"
$clin3BulletChar.Font.NameOther = "Aptos"
$clin3BulletChar.Font.Size = 8

# --- CLIN4 Char: Calibri -> Aptos, 6pt -> 5.5pt (sz 12 -> 11) -------
$clin4Char = $d.Styles("CLIN4 Char")
$clin4Char.Font.NameAscii = "Aptos"
$clin4Char.Font.NameOther = "Aptos"
$clin4Char.Font.Size = 5.5

Write-Output "Applied MHD2-259 template style updates."
